$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is numeric-looking need an explicit Text
# number format first, otherwise Excel auto-converts the assigned
# string into a floating point number (changing cell type from the
# original inline string).
$textCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D15", "D16", "D18", "D20", "D21", "D22", "D25", "D27", "D28", "D31", "D34", "D36", "D37", "D44", "D46", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.778.37'
$ws.Range("E2").Value = '  -0.59%  '

$ws.Range("D3").Value = '2.552.68'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '303.68'
$ws.Range("E5").Value = '  +1.78%  '

$ws.Range("D6").Value = '98.42'
$ws.Range("E6").Value = '  +6.90%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("E9").Value = '  -0.75%  '

$ws.Range("D10").Value = '36.42'
$ws.Range("E10").Value = '  +1.05%  '

$ws.Range("D11").Value = '0.0808'
$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").Value = '0.117'
$ws.Range("E12").Value = '  +7.79%  '

$ws.Range("D13").Value = '7.55'
$ws.Range("E13").Value = '  -2.35%  '

$ws.Range("D14").Value = '2.538.82'
$ws.Range("E14").Value = '  +0.02%  '

$ws.Range("D15").Value = '0.881'
$ws.Range("E15").Value = '  +1.14%  '

$ws.Range("D16").Value = '14.84'
$ws.Range("E16").Value = '  +4.75%  '

$ws.Range("D17").Value = '42.843.36'
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").Value = '13.23'
$ws.Range("E18").Value = '  +5.54%  '

$ws.Range("E19").Value = '  +0.74%  '

$ws.Range("D20").Value = '6.61'
$ws.Range("E20").Value = '  -0.97%  '

$ws.Range("D21").Value = '71.62'
$ws.Range("E21").Value = '  -0.71%  '

$ws.Range("D22").Value = '254.97'
$ws.Range("E22").Value = '  -2.34%  '

$ws.Range("E23").Value = '  +1.99%  '

$ws.Range("E24").Value = '  -2.15%  '

$ws.Range("D25").Value = '27.73'
$ws.Range("E25").Value = '  -6.46%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("D27").Value = '10.08'
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").Value = '38.02'
$ws.Range("E28").Value = '  +3.92%  '

$ws.Range("E29").Value = '  -1.40%  '

$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("D31").Value = '156.72'
$ws.Range("E31").Value = '  +2.22%  '

$ws.Range("E32").Value = '  +0.23%  '

$ws.Range("E33").Value = '  +0.78%  '

$ws.Range("D34").Value = '0.0807'
$ws.Range("E34").Value = '  +1.75%  '

$ws.Range("E35").Value = '  -3.04%  '

$ws.Range("D36").Value = '26.67'
$ws.Range("E36").Value = '  +10.70%  '

$ws.Range("D37").Value = '18.59'
$ws.Range("E37").Value = '  +10.55%  '

$ws.Range("E38").Value = '  +0.51%  '

$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("E40").Value = '  +32.67%  '

$ws.Range("E41").Value = '  -1.65%  '

$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("E43").Value = '  -2.54%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.19%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.062.05'
$ws.Range("E45").Value = '  -0.65%  '

$ws.Range("D46").Value = '87.52'
$ws.Range("E46").Value = '  +2.14%  '

$ws.Range("D47").Value = '9.17'
$ws.Range("E47").Value = '  +5.94%  '

$ws.Range("D48").Value = '2.802.70'
$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("D49").Value = '74.85'
$ws.Range("E49").Value = '  +7.85%  '

$ws.Range("D50").Value = '103.19'
$ws.Range("E50").Value = '  -1.11%  '

$ws.Range("E51").Value = '  +1.31%  '
